$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value. "numeric" values are written with a
# temporary Text number format so Excel does not coerce strings like "94.00" or
# "0.000007815" into doubles (which would silently drop precision/trailing zeros).
# The number format is restored to the default "Normal" style right after, so no
# residual formatting is left behind on the cell.
$updates = @(
    @{ Cell = "D2"; Value = '30.256.30'; Numeric = $false }
    @{ Cell = "E2"; Value = '  -3.57%  '; Numeric = $false }
    @{ Cell = "D3"; Value = '1.927.01'; Numeric = $false }
    @{ Cell = "E3"; Value = '  -3.01%  '; Numeric = $false }
    @{ Cell = "E4"; Value = '  +0.52%  '; Numeric = $false }
    @{ Cell = "D5"; Value = '246.71'; Numeric = $true }
    @{ Cell = "E5"; Value = '  -2.49%  '; Numeric = $false }
    @{ Cell = "D6"; Value = '0.7002'; Numeric = $true }
    @{ Cell = "E6"; Value = '  -11.55%  '; Numeric = $false }
    @{ Cell = "E7"; Value = '  +0.32%  '; Numeric = $false }
    @{ Cell = "D8"; Value = '0.3228'; Numeric = $true }
    @{ Cell = "E8"; Value = '  -5.76%  '; Numeric = $false }
    @{ Cell = "D9"; Value = '26.42'; Numeric = $true }
    @{ Cell = "E9"; Value = '  +2.72%  '; Numeric = $false }
    @{ Cell = "D10"; Value = '0.06818'; Numeric = $true }
    @{ Cell = "E10"; Value = '  -2.06%  '; Numeric = $false }
    @{ Cell = "D11"; Value = '0.7939'; Numeric = $true }
    @{ Cell = "E11"; Value = '  -5.83%  '; Numeric = $false }
    @{ Cell = "D12"; Value = '0.07942'; Numeric = $true }
    @{ Cell = "E12"; Value = '  -2.27%  '; Numeric = $false }
    @{ Cell = "D13"; Value = '1.926.35'; Numeric = $false }
    @{ Cell = "E13"; Value = '  -2.74%  '; Numeric = $false }
    @{ Cell = "D14"; Value = '5.391'; Numeric = $true }
    @{ Cell = "E14"; Value = '  -2.01%  '; Numeric = $false }
    @{ Cell = "D15"; Value = '94.00'; Numeric = $true }
    @{ Cell = "E15"; Value = '  -8.25%  '; Numeric = $false }
    @{ Cell = "E16"; Value = '  +3.57%  '; Numeric = $false }
    @{ Cell = "D17"; Value = '259.44'; Numeric = $true }
    @{ Cell = "E17"; Value = '  -5.77%  '; Numeric = $false }
    @{ Cell = "D18"; Value = '30.261.39'; Numeric = $false }
    @{ Cell = "E18"; Value = '  -3.32%  '; Numeric = $false }
    @{ Cell = "D19"; Value = '5.845'; Numeric = $true }
    @{ Cell = "E19"; Value = '  +2.82%  '; Numeric = $false }
    @{ Cell = "D20"; Value = '0.000007815'; Numeric = $true }
    @{ Cell = "E20"; Value = '  -0.93%  '; Numeric = $false }
    @{ Cell = "D21"; Value = '2.178.49'; Numeric = $false }
    @{ Cell = "E21"; Value = '  -2.57%  '; Numeric = $false }
    @{ Cell = "D22"; Value = '1.001'; Numeric = $true }
    @{ Cell = "E22"; Value = '  +0.43%  '; Numeric = $false }
    @{ Cell = "D23"; Value = '1.002'; Numeric = $true }
    @{ Cell = "E23"; Value = '  +0.50%  '; Numeric = $false }
    @{ Cell = "D24"; Value = '6.810'; Numeric = $true }
    @{ Cell = "E24"; Value = '  -0.75%  '; Numeric = $false }
    @{ Cell = "D25"; Value = '9.604'; Numeric = $true }
    @{ Cell = "E25"; Value = '  -0.49%  '; Numeric = $false }
    @{ Cell = "D26"; Value = '158.56'; Numeric = $true }
    @{ Cell = "E26"; Value = '  -4.80%  '; Numeric = $false }
    @{ Cell = "D27"; Value = '18.77'; Numeric = $true }
    @{ Cell = "E27"; Value = '  -4.38%  '; Numeric = $false }
    @{ Cell = "D28"; Value = '0.1313'; Numeric = $true }
    @{ Cell = "E28"; Value = '  -15.43%  '; Numeric = $false }
    @{ Cell = "D29"; Value = '2.222'; Numeric = $true }
    @{ Cell = "E29"; Value = '  -2.66%  '; Numeric = $false }
    @{ Cell = "B30"; Value = 'PancakeSwap'; Numeric = $false }
    @{ Cell = "C30"; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; Numeric = $false }
    @{ Cell = "D30"; Value = '1.552'; Numeric = $true }
    @{ Cell = "E30"; Value = '  -0.64%  '; Numeric = $false }
    @{ Cell = "B31"; Value = 'Toncoin'; Numeric = $false }
    @{ Cell = "C31"; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; Numeric = $false }
    @{ Cell = "D31"; Value = '1.340'; Numeric = $true }
    @{ Cell = "E31"; Value = '  -0.25%  '; Numeric = $false }
    @{ Cell = "D32"; Value = '4.394'; Numeric = $true }
    @{ Cell = "E32"; Value = '  -3.63%  '; Numeric = $false }
    @{ Cell = "D33"; Value = '4.179'; Numeric = $true }
    @{ Cell = "E33"; Value = '  -3.59%  '; Numeric = $false }
    @{ Cell = "D34"; Value = '0.05036'; Numeric = $true }
    @{ Cell = "E34"; Value = '  -3.29%  '; Numeric = $false }
    @{ Cell = "D35"; Value = '1.191'; Numeric = $true }
    @{ Cell = "E35"; Value = '  -2.16%  '; Numeric = $false }
    @{ Cell = "D36"; Value = '0.7460'; Numeric = $true }
    @{ Cell = "E36"; Value = '  +0.18%  '; Numeric = $false }
    @{ Cell = "D37"; Value = '2.708'; Numeric = $true }
    @{ Cell = "E37"; Value = '  -3.09%  '; Numeric = $false }
    @{ Cell = "D38"; Value = '0.01918'; Numeric = $true }
    @{ Cell = "E38"; Value = '  -3.33%  '; Numeric = $false }
    @{ Cell = "D39"; Value = '2.792'; Numeric = $true }
    @{ Cell = "E39"; Value = '  -3.89%  '; Numeric = $false }
    @{ Cell = "D40"; Value = '80.03'; Numeric = $true }
    @{ Cell = "E40"; Value = '  +1.79%  '; Numeric = $false }
    @{ Cell = "D41"; Value = '6.512'; Numeric = $true }
    @{ Cell = "E41"; Value = '  -1.61%  '; Numeric = $false }
    @{ Cell = "D42"; Value = '2.046'; Numeric = $true }
    @{ Cell = "E42"; Value = '  -1.67%  '; Numeric = $false }
    @{ Cell = "D43"; Value = '0.4398'; Numeric = $true }
    @{ Cell = "E43"; Value = '  -5.93%  '; Numeric = $false }
    @{ Cell = "E44"; Value = '  +0.23%  '; Numeric = $false }
    @{ Cell = "D45"; Value = '0.8331'; Numeric = $true }
    @{ Cell = "E45"; Value = '  -2.53%  '; Numeric = $false }
    @{ Cell = "D46"; Value = '101.61'; Numeric = $true }
    @{ Cell = "E46"; Value = '  -3.94%  '; Numeric = $false }
    @{ Cell = "D47"; Value = '9.732'; Numeric = $true }
    @{ Cell = "E47"; Value = '  -2.67%  '; Numeric = $false }
    @{ Cell = "D48"; Value = '7.198'; Numeric = $true }
    @{ Cell = "E48"; Value = '  -4.04%  '; Numeric = $false }
    @{ Cell = "D49"; Value = '35.76'; Numeric = $true }
    @{ Cell = "E49"; Value = '  -1.97%  '; Numeric = $false }
    @{ Cell = "D50"; Value = '2.811'; Numeric = $true }
    @{ Cell = "E50"; Value = '  +31.03%  '; Numeric = $false }
    @{ Cell = "B51"; Value = 'Cronos'; Numeric = $false }
    @{ Cell = "C51"; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; Numeric = $false }
    @{ Cell = "D51"; Value = '0.05935'; Numeric = $true }
    @{ Cell = "E51"; Value = '  +0.15%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Numeric) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
